# Update the "sandboxed" schema on slide 2: shrink/move the Rounded
# Rectangle (id 79) and the Round Same Side Corner Rectangle (id 80) so the
# "apps" band sits in-between the VEE diagram, highlighting OTA.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Rounded Rectangle 105 (the outer sandbox frame)
$rect = $s.Shapes.Item(15)
$rect.Top = 144.86276
$rect.Height = 144.50458

# Round Same Side Corner Rectangle 106 (the rotated cap/lid)
$cap = $s.Shapes.Item(16)
$cap.Left = 317.88608
$cap.Top = 144.86284
$cap.Height = 28.6306305
